$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet's tab / workbook sheet name
$ws.Name = "UniformF"

# Add new row 16 with Gaussian Quadrature scheme data (HexGrid-60degTilt5degRes), mirroring row 15
$ws.Range("A16").Value2 = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B16").Value2 = $ws.Range("B15").Value2

$ws.Range("C16:M16").Value2 = 1
